$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.647.78'
$ws.Range("E2").Value = '  +1.05%  '
$ws.Range("D3").Value = '1.862.62'
$ws.Range("E3").Value = '  +0.03%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9991'
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '234.45'
$ws.Range("E5").Value = '  +0.14%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9992'
$ws.Range("E6").Value = '  -0.06%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4705'
$ws.Range("E7").Value = '  -1.19%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2756'
$ws.Range("E8").Value = '  +0.16%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06360'
$ws.Range("E9").Value = '  -1.35%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '17.44'
$ws.Range("E10").Value = '  +8.11%  '
$ws.Range("D11").Value = '1.861.73'
$ws.Range("E11").Value = '  +0.50%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07469'
$ws.Range("E12").Value = '  +0.53%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.959'
$ws.Range("E13").Value = '  -0.77%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '84.96'
$ws.Range("E14").Value = '  -1.30%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6298'
$ws.Range("E15").Value = '  -0.65%  '
$ws.Range("D16").Value = '30.593.89'
$ws.Range("E16").Value = '  +0.93%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '240.02'
$ws.Range("E17").Value = '  +3.23%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.9985'
$ws.Range("E18").Value = '  -0.13%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.69'
$ws.Range("E19").Value = '  -1.01%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007359'
$ws.Range("E20").Value = '  -0.49%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9985'
$ws.Range("E21").Value = '  -0.13%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.985'
$ws.Range("E22").Value = '  -2.47%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.956'
$ws.Range("E23").Value = '  -1.09%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.269'
$ws.Range("E24").Value = '  -0.33%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '166.48'
$ws.Range("E25").Value = '  -0.66%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '18.14'
$ws.Range("E26").Value = '  +1.22%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.885'
$ws.Range("E27").Value = '  +1.24%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.1028'
$ws.Range("E28").Value = '  +2.15%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.380'
$ws.Range("E29").Value = '  -0.25%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.102'
$ws.Range("E30").Value = '  -3.18%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.856'
$ws.Range("E31").Value = '  -1.64%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.04929'
$ws.Range("E32").Value = '  +0.53%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.150'
$ws.Range("E33").Value = '  -0.09%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7085'
$ws.Range("E34").Value = '  -2.55%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.697'
$ws.Range("E35").Value = '  +0.14%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.01912'
$ws.Range("E36").Value = '  -1.77%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.692'
$ws.Range("E37").Value = '  +2.26%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.8823'
$ws.Range("E38").Value = '  -2.96%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.988'
$ws.Range("E39").Value = '  -0.25%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '105.94'
$ws.Range("E40").Value = '  +0.29%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9983'
$ws.Range("E41").Value = '  -0.15%  '
$ws.Range("E42").Value = '  -0.62%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.551'
$ws.Range("E43").Value = '  -0.02%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '7.200'
$ws.Range("E44").Value = '  +1.73%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.1232'
$ws.Range("E45").Value = '  +1.90%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '61.73'
$ws.Range("E46").Value = '  +0.40%  '
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.652'
$ws.Range("E47").Value = '  -1.21%  '
$ws.Range("B48").Value = 'Elrond'
$ws.Range("C48").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '33.68'
$ws.Range("E48").Value = '  +1.77%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.05562'
$ws.Range("E49").Value = '  -0.88%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.378'
$ws.Range("E50").Value = '  -1.81%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.3711'
$ws.Range("E51").Value = '  -0.16%  '
